$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '57.941.78'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").Value = '2.352.66'
$ws.Range("E3").Value = '  +1.19%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '543.21'
$ws.Range("E5").Value = '  +0.92%  '
$ws.Range("D6").Value = '134.58'
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("E7").Value = '  +0.68%  '
$ws.Range("D8").Value = '0.557'
$ws.Range("E8").Value = '  +4.21%  '
$ws.Range("D9").Value = '0.102'
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("D10").Value = '5.60'
$ws.Range("E10").Value = '  +3.14%  '
$ws.Range("E11").Value = '  -1.24%  '
$ws.Range("D12").Value = '0.355'
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").Value = '23.82'
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("D14").Value = '2.771.17'
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").Value = '57.914.99'
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").Value = '0.0000134'
$ws.Range("E16").Value = '  +0.54%  '
$ws.Range("D17").Value = '2.370.29'
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("D18").Value = '10.78'
$ws.Range("E18").Value = '  +2.72%  '
$ws.Range("D19").Value = '330.95'
$ws.Range("E19").Value = '  -2.04%  '
$ws.Range("D20").Value = '4.28'
$ws.Range("E20").Value = '  +1.31%  '
$ws.Range("D21").Value = '6.76'
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").Value = '62.47'
$ws.Range("E23").Value = '  +1.05%  '
$ws.Range("E24").Value = '  -1.21%  '
$ws.Range("D25").Value = '8.47'
$ws.Range("E25").Value = '  -1.44%  '
$ws.Range("D26").Value = '0.996'
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").Value = '1.36'
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("D29").Value = '170.84'
$ws.Range("E29").Value = '  -1.56%  '
$ws.Range("D30").Value = '0.0₃0736'
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("D31").Value = '6.15'
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '18.49'
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("B33").Value = 'SuiNetwork'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D33").Value = '1.02'
$ws.Range("E33").Value = '  +2.59%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.84%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '4.22'
$ws.Range("E36").Value = '  +3.82%  '
$ws.Range("D37").Value = '1.24'
$ws.Range("E37").Value = '  -1.74%  '
$ws.Range("D38").Value = '1.63'
$ws.Range("E38").Value = '  +2.25%  '
$ws.Range("D39").Value = '39.33'
$ws.Range("E39").Value = '  +0.22%  '
$ws.Range("D40").Value = '144.99'
$ws.Range("E40").Value = '  -2.30%  '
$ws.Range("D41").Value = '292.96'
$ws.Range("E41").Value = '  +2.73%  '
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("D43").Value = '3.64'
$ws.Range("E43").Value = '  -0.12%  '
$ws.Range("D44").Value = '0.0944'
$ws.Range("E44").Value = '  +1.47%  '
$ws.Range("D45").Value = '19.23'
$ws.Range("E45").Value = '  +1.79%  '
$ws.Range("D46").Value = '0.0506'
$ws.Range("E46").Value = '  +0.71%  '
$ws.Range("D47").Value = '0.564'
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("D48").Value = '0.0224'
$ws.Range("E48").Value = '  +2.81%  '
$ws.Range("D49").Value = '17.54'
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("E50").Value = '  -0.11%  '
$ws.Range("D51").Value = '11.05'
$ws.Range("E51").Value = '  +0.11%  '
